$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '54.432.59'
$ws.Range('E2').Value = '  -6.89%  '

$ws.Range('D3').Value = '2.444.68'
$ws.Range('E3').Value = '  -10.07%  '

$ws.Range('E4').Value = '  +0.13%  '

$r = $ws.Range('D5')
$r.NumberFormat = '@'
$r.Value = '465.67'
$r.Style = 'Normal'
$ws.Range('E5').Value = '  -7.46%  '

$r = $ws.Range('D6')
$r.NumberFormat = '@'
$r.Value = '132.88'
$r.Style = 'Normal'
$ws.Range('E6').Value = '  -5.64%  '

$r = $ws.Range('D7')
$r.NumberFormat = '@'
$r.Value = '0.998'
$r.Style = 'Normal'
$ws.Range('E7').Value = '  +0.09%  '

$r = $ws.Range('D8')
$r.NumberFormat = '@'
$r.Value = '0.492'
$r.Style = 'Normal'
$ws.Range('E8').Value = '  -6.88%  '

$ws.Range('D9').Value = '2.445.17'
$ws.Range('E9').Value = '  -10.39%  '

$r = $ws.Range('D10')
$r.NumberFormat = '@'
$r.Value = '0.0952'
$r.Style = 'Normal'
$ws.Range('E10').Value = '  -9.07%  '

$r = $ws.Range('D11')
$r.NumberFormat = '@'
$r.Value = '5.31'
$r.Style = 'Normal'
$ws.Range('E11').Value = '  -12.07%  '

$r = $ws.Range('D12')
$r.NumberFormat = '@'
$r.Value = '0.316'
$r.Style = 'Normal'
$ws.Range('E12').Value = '  -8.97%  '

$r = $ws.Range('D13')
$r.NumberFormat = '@'
$r.Value = '0.122'
$r.Style = 'Normal'
$ws.Range('E13').Value = '  -3.80%  '

$ws.Range('D14').Value = '2.868.05'
$ws.Range('E14').Value = '  -10.35%  '

$ws.Range('D15').Value = '54.525.18'
$ws.Range('E15').Value = '  -6.89%  '

$ws.Range('E16').Value = '  -1.48%  '

$r = $ws.Range('D17')
$r.NumberFormat = '@'
$r.Value = '19.78'
$r.Style = 'Normal'
$ws.Range('E17').Value = '  -8.73%  '

$ws.Range('D18').Value = '2.435.93'
$ws.Range('E18').Value = '  -10.83%  '

$r = $ws.Range('D19')
$r.NumberFormat = '@'
$r.Value = '4.19'
$r.Style = 'Normal'
$ws.Range('E19').Value = '  -11.83%  '

$r = $ws.Range('D20')
$r.NumberFormat = '@'
$r.Value = '310.45'
$r.Style = 'Normal'
$ws.Range('E20').Value = '  -9.53%  '

$r = $ws.Range('D21')
$r.NumberFormat = '@'
$r.Value = '9.54'
$r.Style = 'Normal'
$ws.Range('E21').Value = '  -12.85%  '

$r = $ws.Range('D22')
$r.NumberFormat = '@'
$r.Value = '0.999'
$r.Style = 'Normal'
$ws.Range('E22').Value = '  +0.27%  '

$r = $ws.Range('D23')
$r.NumberFormat = '@'
$r.Value = '5.69'
$r.Style = 'Normal'
$ws.Range('E23').Value = '  +0.47%  '

$r = $ws.Range('D24')
$r.NumberFormat = '@'
$r.Value = '5.39'
$r.Style = 'Normal'
$ws.Range('E24').Value = '  -13.82%  '

$r = $ws.Range('D25')
$r.NumberFormat = '@'
$r.Value = '56.30'
$r.Style = 'Normal'
$ws.Range('E25').Value = '  -10.29%  '

$ws.Range('E26').Value = '  +0.75%  '

$ws.Range('E27').Value = '  -9.59%  '

$ws.Range('E28').Value = '  -9.95%  '

$ws.Range('D29').Value = '2.528.52'
$ws.Range('E29').Value = '  -11.09%  '

$ws.Range('E30').Value = '  -4.96%  '

$r = $ws.Range('D31')
$r.NumberFormat = '@'
$r.Value = '1.00'
$r.Style = 'Normal'
$ws.Range('E31').Value = '  +0.10%  '

$ws.Range('D32').Value = '0.0₃0716'
$ws.Range('E32').Value = '  -13.55%  '

$r = $ws.Range('D33')
$r.NumberFormat = '@'
$r.Value = '146.92'
$r.Style = 'Normal'
$ws.Range('E33').Value = '  -2.71%  '

$ws.Range('E34').Value = '  -7.13%  '

$r = $ws.Range('D35')
$r.NumberFormat = '@'
$r.Value = '1.43'
$r.Style = 'Normal'
$ws.Range('E35').Value = '  -10.52%  '

$r = $ws.Range('D36')
$r.NumberFormat = '@'
$r.Value = '5.00'
$r.Style = 'Normal'
$ws.Range('E36').Value = '  -7.74%  '

$ws.Range('E37').Value = '  -15.22%  '

$ws.Range('E38').Value = '  -6.87%  '

$r = $ws.Range('D39')
$r.NumberFormat = '@'
$r.Value = '0.805'
$r.Style = 'Normal'
$ws.Range('E39').Value = '  -15.17%  '

$r = $ws.Range('D40')
$r.NumberFormat = '@'
$r.Value = '0.997'
$r.Style = 'Normal'
$ws.Range('E40').Value = '  +0.16%  '

$r = $ws.Range('D41')
$r.NumberFormat = '@'
$r.Value = '32.88'
$r.Style = 'Normal'
$ws.Range('E41').Value = '  -8.32%  '

$r = $ws.Range('D42')
$r.NumberFormat = '@'
$r.Value = '0.597'
$r.Style = 'Normal'
$ws.Range('E42').Value = '  -0.62%  '

$ws.Range('E43').Value = '  -6.25%  '

$ws.Range('E44').Value = '  -8.51%  '

$ws.Range('B45').Value = 'WhiteBITCoin'
$ws.Range('C45').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$r = $ws.Range('D45')
$r.NumberFormat = '@'
$r.Value = '10.07'
$r.Style = 'Normal'
$ws.Range('E45').Value = '  -2.92%  '

$ws.Range('B46').Value = 'Stacks'
$ws.Range('C46').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$r = $ws.Range('D46')
$r.NumberFormat = '@'
$r.Value = '1.24'
$r.Style = 'Normal'
$ws.Range('E46').Value = '  -11.14%  '

$ws.Range('D47').Value = '1.934.82'
$ws.Range('E47').Value = '  -11.49%  '

$ws.Range('E48').Value = '  -0.12%  '

$ws.Range('E49').Value = '  -4.17%  '

$ws.Range('B50').Value = 'EnergySwap'
$ws.Range('C50').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$r = $ws.Range('D50')
$r.NumberFormat = '@'
$r.Value = '16.53'
$r.Style = 'Normal'
$ws.Range('E50').Value = '  -12.68%  '

$ws.Range('B51').Value = 'Bittensor'
$ws.Range('C51').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$r = $ws.Range('D51')
$r.NumberFormat = '@'
$r.Value = '229.42'
$r.Style = 'Normal'
$ws.Range('E51').Value = '  +1.31%  '
